$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 13334357
$ws.Cells.Item(6, 10).Value = 1599.1333
$ws.Cells.Item(6, 12).Value = 4797.3999
$ws.Cells.Item(6, 14).Value = -5021.3999

$ws.Cells.Item(8, 8).Value = 23809722
$ws.Cells.Item(8, 9).Value = 55555620
$ws.Cells.Item(8, 10).Value = 299
$ws.Cells.Item(8, 11).Value = 166666860
$ws.Cells.Item(8, 12).Value = 897
$ws.Cells.Item(8, 13).Value = -166666721
$ws.Cells.Item(8, 14).Value = -1175

$ws.Cells.Item(18, 8).Value = 7181.7
$ws.Cells.Item(18, 9).Value = 7868.5557
$ws.Cells.Item(18, 11).Value = 7868.5557
$ws.Cells.Item(18, 13).Value = -7584.5557

$ws.Cells.Item(19, 8).Value = 430.25
$ws.Cells.Item(19, 9).Value = 176
$ws.Cells.Item(19, 10).Value = 515
$ws.Cells.Item(19, 11).Value = 176
$ws.Cells.Item(19, 12).Value = 515
$ws.Cells.Item(19, 13).Value = -1
$ws.Cells.Item(19, 14).Value = -865

$ws.Cells.Item(76, 8).Value = 4642.9355
$ws.Cells.Item(76, 9).Value = 3869.375
$ws.Cells.Item(76, 10).Value = 5468.067
$ws.Cells.Item(76, 11).Value = 3869.375
$ws.Cells.Item(76, 12).Value = 5468.067
$ws.Cells.Item(76, 13).Value = -3554.375
$ws.Cells.Item(76, 14).Value = -6098.067

$ws.Cells.Item(79, 8).Value = 4642.9355
$ws.Cells.Item(79, 9).Value = 3869.375
$ws.Cells.Item(79, 10).Value = 5468.067
$ws.Cells.Item(79, 11).Value = 3869.375
$ws.Cells.Item(79, 12).Value = 5468.067
$ws.Cells.Item(79, 13).Value = -2777.375
$ws.Cells.Item(79, 14).Value = -7652.067

$ws.Cells.Item(80, 8).Value = 50937.586
$ws.Cells.Item(80, 9).Value = 81083.28
$ws.Cells.Item(80, 11).Value = 243249.84
$ws.Cells.Item(80, 13).Value = -242251.84

$ws.Cells.Item(83, 8).Value = 50937.586
$ws.Cells.Item(83, 9).Value = 81083.28
$ws.Cells.Item(83, 11).Value = 729749.52
$ws.Cells.Item(83, 13).Value = -724757.52

$ws.Cells.Item(112, 8).Value = 2161.7646
$ws.Cells.Item(112, 9).Value = 2959.8
$ws.Cells.Item(112, 10).Value = 1829.25
$ws.Cells.Item(112, 11).Value = 8879.400000000001
$ws.Cells.Item(112, 12).Value = 5487.75
$ws.Cells.Item(112, 13).Value = -7771.400000000001
$ws.Cells.Item(112, 14).Value = -7703.75

$ws.Cells.Item(132, 8).Value = 2329458.8
$ws.Cells.Item(132, 9).Value = 3938.2307
$ws.Cells.Item(132, 11).Value = 11814.6921
$ws.Cells.Item(132, 13).Value = -9284.6921

$ws.Cells.Item(137, 8).Value = 6539.964
$ws.Cells.Item(137, 9).Value = 8829.5
$ws.Cells.Item(137, 10).Value = 2201.8948
$ws.Cells.Item(137, 11).Value = 26488.5
$ws.Cells.Item(137, 12).Value = 6605.6844
$ws.Cells.Item(137, 13).Value = -23938.5
$ws.Cells.Item(137, 14).Value = -11705.6844

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6754.6616
$ws.Cells.Item(32, 9).Value = 6601.746
$ws.Cells.Item(32, 11).Value = 6601.746
$ws.Cells.Item(32, 13).Value = -6314.746

$ws.Cells.Item(45, 8).Value = 6305.5
$ws.Cells.Item(45, 9).Value = 5694.4443
$ws.Cells.Item(45, 10).Value = 7405.4
$ws.Cells.Item(45, 11).Value = 5694.4443
$ws.Cells.Item(45, 12).Value = 7405.4
$ws.Cells.Item(45, 13).Value = -5317.4443
$ws.Cells.Item(45, 14).Value = -8159.4

$ws.Cells.Item(74, 8).Value = 2226.3125
$ws.Cells.Item(74, 9).Value = 1020.1667
$ws.Cells.Item(74, 11).Value = 1020.1667
$ws.Cells.Item(74, 13).Value = -146.1667

$ws.Cells.Item(77, 8).Value = 2226.3125
$ws.Cells.Item(77, 9).Value = 1020.1667
$ws.Cells.Item(77, 11).Value = 5100.8335
$ws.Cells.Item(77, 13).Value = -732.8334999999997

$ws.Cells.Item(102, 8).Value = 12905.629
$ws.Cells.Item(102, 9).Value = 22066.934
$ws.Cells.Item(102, 11).Value = 22066.934
$ws.Cells.Item(102, 13).Value = -20444.934

$ws.Cells.Item(132, 8).Value = 5474.385
$ws.Cells.Item(132, 9).Value = 1785
$ws.Cells.Item(132, 11).Value = 5355
$ws.Cells.Item(132, 13).Value = -2825

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 7374.0454
$ws.Cells.Item(105, 9).Value = 8196
$ws.Cells.Item(105, 11).Value = 8196
$ws.Cells.Item(105, 13).Value = -6449

$ws.Cells.Item(107, 8).Value = 4498.1665
$ws.Cells.Item(107, 9).Value = 5247.5
$ws.Cells.Item(107, 10).Value = 2999.5
$ws.Cells.Item(107, 11).Value = 5247.5
$ws.Cells.Item(107, 12).Value = 2999.5
$ws.Cells.Item(107, 13).Value = -3327.5
$ws.Cells.Item(107, 14).Value = -6839.5

$ws.Cells.Item(134, 8).Value = 5749.5586
$ws.Cells.Item(134, 9).Value = 6982.7144
$ws.Cells.Item(134, 11).Value = 20948.1432
$ws.Cells.Item(134, 13).Value = -18413.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 13846.571
$ws.Cells.Item(31, 9).Value = 15715.546
$ws.Cells.Item(31, 10).Value = 6993.6665
$ws.Cells.Item(31, 11).Value = 15715.546
$ws.Cells.Item(31, 12).Value = 6993.6665
$ws.Cells.Item(31, 13).Value = -15420.546
$ws.Cells.Item(31, 14).Value = -7583.6665

$ws.Cells.Item(34, 8).Value = 13846.571
$ws.Cells.Item(34, 9).Value = 15715.546
$ws.Cells.Item(34, 10).Value = 6993.6665
$ws.Cells.Item(34, 11).Value = 15715.546
$ws.Cells.Item(34, 12).Value = 6993.6665
$ws.Cells.Item(34, 13).Value = -15513.546
$ws.Cells.Item(34, 14).Value = -7397.6665

$ws.Cells.Item(134, 8).Value = 6355.9287
$ws.Cells.Item(134, 9).Value = 5367.7144
$ws.Cells.Item(134, 11).Value = 16103.1432
$ws.Cells.Item(134, 13).Value = -13568.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(17, 8).Value = 1827.375
$ws.Cells.Item(17, 10).Value = 5995
$ws.Cells.Item(17, 12).Value = 17985
$ws.Cells.Item(17, 14).Value = -18323

$ws.Cells.Item(42, 8).Value = 7500
$ws.Cells.Item(42, 10).Value = 7500
$ws.Cells.Item(42, 12).Value = 22500
$ws.Cells.Item(42, 14).Value = -23568

$ws.Cells.Item(128, 8).Value = 549999.5
$ws.Cells.Item(128, 9).Value = 549999.5
$ws.Cells.Item(128, 11).Value = 1649998.5
$ws.Cells.Item(128, 13).Value = -1645018.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 5965.613
$ws.Cells.Item(97, 9).Value = 6371.44
$ws.Cells.Item(97, 11).Value = 6371.44
$ws.Cells.Item(97, 13).Value = -5875.44

$ws.Cells.Item(132, 8).Value = 2426.1
$ws.Cells.Item(132, 9).Value = 2441.862
$ws.Cells.Item(132, 11).Value = 7325.586
$ws.Cells.Item(132, 13).Value = -4795.586

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 22399.6
$ws.Cells.Item(7, 10).Value = 7749.75
$ws.Cells.Item(7, 12).Value = 7749.75
$ws.Cells.Item(7, 14).Value = -7973.75

$ws.Cells.Item(126, 8).Value = 22399.6
$ws.Cells.Item(126, 10).Value = 7749.75
$ws.Cells.Item(126, 12).Value = 23249.25
$ws.Cells.Item(126, 14).Value = -28189.25

$ws.Cells.Item(136, 8).Value = 4999.184
$ws.Cells.Item(136, 9).Value = 3775.2856
$ws.Cells.Item(136, 10).Value = 5713.125
$ws.Cells.Item(136, 11).Value = 11325.8568
$ws.Cells.Item(136, 12).Value = 17139.375
$ws.Cells.Item(136, 13).Value = -8775.856800000001
$ws.Cells.Item(136, 14).Value = -22239.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(75, 8).Value = 39583.168
$ws.Cells.Item(75, 9).Value = 27500
$ws.Cells.Item(75, 11).Value = 27500
$ws.Cells.Item(75, 13).Value = -26564

$ws.Cells.Item(78, 8).Value = 39583.168
$ws.Cells.Item(78, 9).Value = 27500
$ws.Cells.Item(78, 11).Value = 82500
$ws.Cells.Item(78, 13).Value = -77820

$ws.Cells.Item(81, 8).Value = 26090.25
$ws.Cells.Item(81, 9).Value = 50550
$ws.Cells.Item(81, 11).Value = 101100
$ws.Cells.Item(81, 13).Value = -100039

$ws.Cells.Item(84, 8).Value = 26090.25
$ws.Cells.Item(84, 9).Value = 50550
$ws.Cells.Item(84, 11).Value = 505500
$ws.Cells.Item(84, 13).Value = -500196

$ws.Cells.Item(132, 8).Value = 10913.155
$ws.Cells.Item(132, 9).Value = 13088.226
$ws.Cells.Item(132, 11).Value = 39264.678
$ws.Cells.Item(132, 13).Value = -36734.678

$ws.Cells.Item(136, 8).Value = 396349.38
$ws.Cells.Item(136, 9).Value = 441277.84
$ws.Cells.Item(136, 10).Value = 3225.5
$ws.Cells.Item(136, 11).Value = 1323833.52
$ws.Cells.Item(136, 12).Value = 9676.5
$ws.Cells.Item(136, 13).Value = -14776.5
